$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "  Jacques Offenbach (14) is enrolled in the Paris Conservatoire."
#      -> "  Jacob (Jacques) Offenbach (14) is enrolled in the Paris Conservatoire."
# split across 5 runs: "  " | "Jacob (" | "Jacques" | ")" | " Offenbach (14)..."
# ---------------------------------------------------------------------------
$r1 = $d.Content
$f1 = $r1.Find
$f1.Text = "Jacques Offenbach"
$f1.Execute() | Out-Null
$jStart = $r1.Start
$jEnd = $jStart + 7   # length of "Jacques"

# Insert the new text fragments first (formatting unchanged so Word's usual
# run-coalescing keeps things simple while we still know exact offsets).
$afterJacques = $d.Range($jEnd, $jEnd)
$afterJacques.InsertBefore(")")

$beforeJacques = $d.Range($jStart, $jStart)
$beforeJacques.InsertBefore("Jacob (")

# Recompute the (now shifted) offsets of each of the 5 pieces.
$prefixStart = $jStart                 # "Jacob ("
$prefixEnd   = $prefixStart + 7
$jacquesStart = $prefixEnd             # "Jacques"
$jacquesEnd   = $jacquesStart + 7
$closeStart   = $jacquesEnd            # ")"
$closeEnd     = $closeStart + 1

# Now force run boundaries between the pieces by toggling formatting on/off.
# (Must be done only after ALL InsertBefore calls in this paragraph, since a
# subsequent text-insertion re-coalesces runs that share identical rPr.)
$prefixRng = $d.Range($prefixStart, $prefixEnd)
$prefixRng.Bold = $true
$prefixRng.Bold = $false

$jacquesRng = $d.Range($jacquesStart, $jacquesEnd)
$jacquesRng.Bold = $true
$jacquesRng.Bold = $false

$closeRng = $d.Range($closeStart, $closeEnd)
$closeRng.Bold = $true
$closeRng.Bold = $false

# ---------------------------------------------------------------------------
# Edit 2: "April 2016" -> "June 2016", split as two runs "June" | " 2016"
# ---------------------------------------------------------------------------
$r2 = $d.Content
$f2 = $r2.Find
$f2.Text = "April 2016"
$f2.Execute() | Out-Null
$aStart = $r2.Start

$aprilRng = $d.Range($aStart, $aStart + 5)   # "April"
$aprilRng.Text = "June"

$juneRng = $d.Range($aStart, $aStart + 4)    # "June" (now 4 chars)
$juneRng.Bold = $true
$juneRng.Bold = $false

Write-Output "edits applied"
